$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 11.83782777605671
$ws.Cells.Item(2, 3).Value = 12.08617999681385
$ws.Cells.Item(2, 4).Value = 5.176933882879548
$ws.Cells.Item(2, 6).Value = 24.36836948232451
$ws.Cells.Item(2, 7).Value = 28.86941921813986
$ws.Cells.Item(2, 8).Value = 14.41038425423528
$ws.Cells.Item(2, 11).Value = 8.253370248803121
$ws.Cells.Item(2, 12).Value = 10.72012794738739
$ws.Cells.Item(2, 15).Value = 21.93241207309113

$ws.Cells.Item(3, 2).Value = 11.52406240072966
$ws.Cells.Item(3, 3).Value = 12.13715618049211
$ws.Cells.Item(3, 4).Value = 5.105982635988368
$ws.Cells.Item(3, 6).Value = 24.4244233285394
$ws.Cells.Item(3, 7).Value = 28.96703019717387
$ws.Cells.Item(3, 8).Value = 14.46015783104575
$ws.Cells.Item(3, 11).Value = 7.973659996703547
$ws.Cells.Item(3, 12).Value = 10.69299948657499
$ws.Cells.Item(3, 15).Value = 22.01606865379292

$ws.Cells.Item(4, 2).Value = 11.32866115244245
$ws.Cells.Item(4, 3).Value = 12.1700509872187
$ws.Cells.Item(4, 4).Value = 5.061309320153142
$ws.Cells.Item(4, 6).Value = 24.46520533753716
$ws.Cells.Item(4, 7).Value = 29.03609647572991
$ws.Cells.Item(4, 8).Value = 14.49294187366782
$ws.Cells.Item(4, 11).Value = 7.795377433984188
$ws.Cells.Item(4, 12).Value = 10.67847097180289
$ws.Cells.Item(4, 15).Value = 22.07200696211601

$ws.Cells.Item(5, 2).Value = 11.24846112425129
$ws.Cells.Item(5, 3).Value = 12.18385821525728
$ws.Cells.Item(5, 4).Value = 5.042836237016681
$ws.Cells.Item(5, 6).Value = 24.4834209588389
$ws.Cells.Item(5, 7).Value = 29.06652808122892
$ws.Cells.Item(5, 8).Value = 14.50686077576121
$ws.Cells.Item(5, 11).Value = 7.721150566228078
$ws.Cells.Item(5, 12).Value = 10.67309018205215
$ws.Cells.Item(5, 15).Value = 22.09595042626428

$ws.Cells.Item(6, 2).Value = 11.2351129889369
$ws.Cells.Item(6, 3).Value = 12.18617522973291
$ws.Cells.Item(6, 4).Value = 5.039752901531255
$ws.Cells.Item(6, 6).Value = 24.48654195960006
$ws.Cells.Item(6, 7).Value = 29.0717190242787
$ws.Cells.Item(6, 8).Value = 14.50920577493422
$ws.Cells.Item(6, 11).Value = 7.708732274133546
$ws.Cells.Item(6, 12).Value = 10.67222942332775
$ws.Cells.Item(6, 15).Value = 22.09999551469836

$ws.Cells.Item(7, 2).Value = 11.32758170313059
$ws.Cells.Item(7, 3).Value = 12.17023556572891
$ws.Cells.Item(7, 4).Value = 5.061061258104378
$ws.Cells.Item(7, 6).Value = 24.46544454053945
$ws.Cells.Item(7, 7).Value = 29.03649764253345
$ws.Cells.Item(7, 8).Value = 14.49312732467573
$ws.Cells.Item(7, 11).Value = 7.794382665970441
$ws.Cells.Item(7, 12).Value = 10.67839621383267
$ws.Cells.Item(7, 15).Value = 22.07232522554088

$ws.Cells.Item(8, 2).Value = 11.73028092802912
$ws.Cells.Item(8, 3).Value = 12.10342621718615
$ws.Cells.Item(8, 4).Value = 5.152706264539614
$ws.Cells.Item(8, 6).Value = 24.38637367882902
$ws.Cells.Item(8, 7).Value = 28.90117417143065
$ws.Cells.Item(8, 8).Value = 14.42708493688178
$ws.Cells.Item(8, 11).Value = 8.158322565485268
$ws.Cells.Item(8, 12).Value = 10.71033512541666
$ws.Cells.Item(8, 15).Value = 21.96030674413357

$ws.Cells.Item(9, 2).Value = 12.49330424759333
$ws.Cells.Item(9, 3).Value = 11.98501585129152
$ws.Cells.Item(9, 4).Value = 5.323178966594573
$ws.Cells.Item(9, 6).Value = 24.28196969146821
$ws.Cells.Item(9, 7).Value = 28.70868782068685
$ws.Cells.Item(9, 8).Value = 14.31520522510352
$ws.Cells.Item(9, 11).Value = 8.817453656491622
$ws.Cells.Item(9, 12).Value = 10.78963724621959
$ws.Cells.Item(9, 15).Value = 21.77699821614117

$ws.Cells.Item(10, 2).Value = 13.03164172259596
$ws.Cells.Item(10, 3).Value = 11.90562547715355
$ws.Cells.Item(10, 4).Value = 5.442223379666257
$ws.Cells.Item(10, 6).Value = 24.23632356047233
$ws.Cells.Item(10, 7).Value = 28.61222929433501
$ws.Cells.Item(10, 8).Value = 14.24374296870409
$ws.Cells.Item(10, 11).Value = 9.265367610914801
$ws.Cells.Item(10, 12).Value = 10.85773942811688
$ws.Cells.Item(10, 15).Value = 21.66458933425186

$ws.Cells.Item(11, 2).Value = 13.27059664719521
$ws.Cells.Item(11, 3).Value = 11.87114386467888
$ws.Cells.Item(11, 4).Value = 5.494920540270128
$ws.Cells.Item(11, 6).Value = 24.22232969536117
$ws.Cells.Item(11, 7).Value = 28.57820900013864
$ws.Cells.Item(11, 8).Value = 14.21356188196343
$ws.Cells.Item(11, 11).Value = 9.460712878248449
$ws.Cells.Item(11, 12).Value = 10.89077847820771
$ws.Cells.Item(11, 15).Value = 21.61830821891612

$ws.Cells.Item(12, 2).Value = 13.36014620560444
$ws.Cells.Item(12, 3).Value = 11.85832024216455
$ws.Cells.Item(12, 4).Value = 5.514657163121803
$ws.Cells.Item(12, 6).Value = 24.21800548403882
$ws.Cells.Item(12, 7).Value = 28.56675091213406
$ws.Cells.Item(12, 8).Value = 14.20246766781572
$ws.Cells.Item(12, 11).Value = 9.533438295114202
$ws.Cells.Item(12, 12).Value = 10.90357820710801
$ws.Cells.Item(12, 15).Value = 21.60148255291973

$ws.Cells.Item(13, 2).Value = 13.34090305329438
$ws.Cells.Item(13, 3).Value = 11.86107165669083
$ws.Cells.Item(13, 4).Value = 5.510416394143176
$ws.Cells.Item(13, 6).Value = 24.21889340391384
$ws.Cells.Item(13, 7).Value = 28.56915515338066
$ws.Cells.Item(13, 8).Value = 14.20484211862124
$ws.Cells.Item(13, 11).Value = 9.517831612464381
$ws.Cells.Item(13, 12).Value = 10.90080884713701
$ws.Cells.Item(13, 15).Value = 21.60507509840261

$ws.Cells.Item(14, 2).Value = 13.27798318441142
$ws.Cells.Item(14, 3).Value = 11.87008417763579
$ws.Cells.Item(14, 4).Value = 5.496548717148602
$ws.Cells.Item(14, 6).Value = 24.2219543947454
$ws.Cells.Item(14, 7).Value = 28.57723774750488
$ws.Cells.Item(14, 8).Value = 14.21264244606865
$ws.Cells.Item(14, 11).Value = 9.466721244849071
$ws.Cells.Item(14, 12).Value = 10.89182578110015
$ws.Cells.Item(14, 15).Value = 21.61690992143075

$ws.Cells.Item(15, 2).Value = 13.2393185048289
$ws.Cells.Item(15, 3).Value = 11.87563502496008
$ws.Cells.Item(15, 4).Value = 5.488025625267431
$ws.Cells.Item(15, 6).Value = 24.22395633794403
$ws.Cells.Item(15, 7).Value = 28.58237429302565
$ws.Cells.Item(15, 8).Value = 14.21746396161864
$ws.Cells.Item(15, 11).Value = 9.435251157106677
$ws.Cells.Item(15, 12).Value = 10.8863607418608
$ws.Cells.Item(15, 15).Value = 21.62425030850136

$ws.Cells.Item(16, 2).Value = 13.01589865265371
$ws.Cells.Item(16, 3).Value = 11.90791170901725
$ws.Cells.Item(16, 4).Value = 5.438749389249564
$ws.Cells.Item(16, 6).Value = 24.23737446558354
$ws.Cells.Item(16, 7).Value = 28.61465160381037
$ws.Cells.Item(16, 8).Value = 14.24576220980931
$ws.Cells.Item(16, 11).Value = 9.252428700147645
$ws.Cells.Item(16, 12).Value = 10.85562104762785
$ws.Cells.Item(16, 15).Value = 21.66771176188544

$ws.Cells.Item(17, 2).Value = 12.87725512653452
$ws.Cells.Item(17, 3).Value = 11.92813001096948
$ws.Cells.Item(17, 4).Value = 5.408140248101735
$ws.Cells.Item(17, 6).Value = 24.24734135169597
$ws.Cells.Item(17, 7).Value = 28.6369830576973
$ws.Cells.Item(17, 8).Value = 14.26371839886867
$ws.Cells.Item(17, 11).Value = 9.13808975895585
$ws.Cells.Item(17, 12).Value = 10.83728538757053
$ws.Cells.Item(17, 15).Value = 21.69561864961958

$ws.Cells.Item(18, 2).Value = 12.79695668744145
$ws.Cells.Item(18, 3).Value = 11.93991285665435
$ws.Cells.Item(18, 4).Value = 5.390398086029889
$ws.Cells.Item(18, 6).Value = 24.25371129792537
$ws.Cells.Item(18, 7).Value = 28.65075515097488
$ws.Cells.Item(18, 8).Value = 14.27426539271283
$ws.Cells.Item(18, 11).Value = 9.071535074988768
$ws.Cells.Item(18, 12).Value = 10.82693348463129
$ws.Cells.Item(18, 15).Value = 21.7121267442305

$ws.Cells.Item(19, 2).Value = 12.7696764795349
$ws.Cells.Item(19, 3).Value = 11.94392877624493
$ws.Cells.Item(19, 4).Value = 5.384367721144194
$ws.Cells.Item(19, 6).Value = 24.25597745135626
$ws.Cells.Item(19, 7).Value = 28.65557724173281
$ws.Cells.Item(19, 8).Value = 14.27787404831966
$ws.Cells.Item(19, 11).Value = 9.048866362415588
$ws.Cells.Item(19, 12).Value = 10.82346209892395
$ws.Cells.Item(19, 15).Value = 21.71779449544519

$ws.Cells.Item(20, 2).Value = 12.89207196984803
$ws.Cells.Item(20, 3).Value = 11.92596182694442
$ws.Cells.Item(20, 4).Value = 5.411412848142628
$ws.Cells.Item(20, 6).Value = 24.24621439438769
$ws.Cells.Item(20, 7).Value = 28.63450977700917
$ws.Cells.Item(20, 8).Value = 14.26178425893953
$ws.Cells.Item(20, 11).Value = 9.150343380862044
$ws.Cells.Item(20, 12).Value = 10.83921719542325
$ws.Cells.Item(20, 15).Value = 21.69260062333158

$ws.Cells.Item(21, 2).Value = 13.29649034774576
$ws.Cells.Item(21, 3).Value = 11.86743064458406
$ws.Cells.Item(21, 4).Value = 5.50062799504441
$ws.Cells.Item(21, 6).Value = 24.22102883976471
$ws.Cells.Item(21, 7).Value = 28.57482497800742
$ws.Cells.Item(21, 8).Value = 14.21034221669475
$ws.Cells.Item(21, 11).Value = 9.481767735106573
$ws.Cells.Item(21, 12).Value = 10.89445655606103
$ws.Cells.Item(21, 15).Value = 21.61341473422783

$ws.Cells.Item(22, 2).Value = 13.55530319746822
$ws.Cells.Item(22, 3).Value = 11.83053949997954
$ws.Cells.Item(22, 4).Value = 5.557656852424416
$ws.Cells.Item(22, 6).Value = 24.21025163245028
$ws.Cells.Item(22, 7).Value = 28.54412410265305
$ws.Cells.Item(22, 8).Value = 14.1786728219045
$ws.Cells.Item(22, 11).Value = 9.691086600554433
$ws.Cells.Item(22, 12).Value = 10.93223682259098
$ws.Cells.Item(22, 15).Value = 21.56574311517723

$ws.Cells.Item(23, 2).Value = 13.41769868191273
$ws.Cells.Item(23, 3).Value = 11.85010469366751
$ws.Cells.Item(23, 4).Value = 5.527339369923474
$ws.Cells.Item(23, 6).Value = 24.21548335311429
$ws.Cells.Item(23, 7).Value = 28.55974762426185
$ws.Cells.Item(23, 8).Value = 14.19539685154125
$ws.Cells.Item(23, 11).Value = 9.580047103389578
$ws.Cells.Item(23, 12).Value = 10.91192182701042
$ws.Cells.Item(23, 15).Value = 21.59081231673591

$ws.Cells.Item(24, 2).Value = 12.88537510760506
$ws.Cells.Item(24, 3).Value = 11.92694156734886
$ws.Cells.Item(24, 4).Value = 5.409933755232972
$ws.Cells.Item(24, 6).Value = 24.24672189880852
$ws.Cells.Item(24, 7).Value = 28.63562504013139
$ws.Cells.Item(24, 8).Value = 14.26265798665592
$ws.Cells.Item(24, 11).Value = 9.144806068698717
$ws.Cells.Item(24, 12).Value = 10.83834323422251
$ws.Cells.Item(24, 15).Value = 21.69396362766147

$ws.Cells.Item(25, 2).Value = 12.29039003554096
$ws.Cells.Item(25, 3).Value = 12.01570781733036
$ws.Cells.Item(25, 4).Value = 5.278111104856812
$ws.Cells.Item(25, 6).Value = 24.30477024989518
$ws.Cells.Item(25, 7).Value = 28.75290032218764
$ws.Cells.Item(25, 8).Value = 14.34358554649128
$ws.Cells.Item(25, 11).Value = 8.645325675760516
$ws.Cells.Item(25, 12).Value = 10.76643074911256
$ws.Cells.Item(25, 15).Value = 21.82268481889774
